$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.597.63'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.460.30'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range('E10').Value = '  -4.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.67'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '2.838.80'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.93%  '
$ws.Range('D16').Value = '2.465.97'
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '41.495.01'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('E20').Value = '  -4.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('E22').Value = '  -2.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  -1.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.89'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.90'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('E42').Value = '  -7.26%  '
$ws.Range('D43').Value = '1.995.65'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0280'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.65%  '
$ws.Range('D48').Value = '2.717.70'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '66.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.53%  '
